# Add a new worksheet "p1.xlsx" at the end of the workbook, mirroring the
# layout used by the existing "pr10.xlsx" / "p11.xlsx" sheets, and populate
# it with the PVRP short-solution data for p1.

$wb = $excel.ActiveWorkbook

# Create the new sheet after the last existing sheet so it lands at the end
# of the tab strip (Worksheets.Add() alone inserts before the active sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "p1.xlsx"

# --- Header / metadata block (rows 1-3) ---------------------------------
$ws.Range("A1").Value = "Filename"
$ws.Range("B1").Value = "p1.xlsx"

$ws.Range("A2").Value = "Selection Type"
$ws.Range("B2").Value = "edu.sru.thangiah.zeus.tr.TRSolutionHierarchy.Heuristics.Selection.TRSmallestAngleClosestDistanceToDepot"

$ws.Range("A3").Value = "Insertion Type"
$ws.Range("B3").Value = "edu.sru.thangiah.zeus.tr.TRSolutionHierarchy.Heuristics.Insertion.TRGreedyInsertion"

# rows 4-5 intentionally left blank

# --- Column headers for the route table (row 6) --------------------------
$ws.Range("A6").Value = "Depot"
$ws.Range("B6").Value = "Truck"
$ws.Range("C6").Value = "Day"
$ws.Range("D6").Value = "Demand"
$ws.Range("E6").Value = "Nodes -->"

# --- Route data rows (7-12) -----------------------------------------------
$data = @{
    7  = @(0.0, 0.0, 0.0, 160.0, -1.0, 51.0, -1.0)
    8  = @(0.0, 0.0, 1.0, 274.0, -1.0, 1.0, 11.0, 12.0, 24.0, 7.0, 26.0, 13.0, 19.0, 23.0, 42.0, 5.0, 17.0, 18.0, 33.0, 41.0, 40.0, -1.0)
    9  = @(0.0, 1.0, 2.0, 0.0, -1.0, -1.0)
    10 = @(0.0, 1.0, 3.0, 259.0, -1.0, 32.0, 6.0, 47.0, 4.0, 37.0, 44.0, 49.0, 10.0, 25.0, 31.0, 22.0, 46.0, 3.0, 48.0, 9.0, 29.0, 21.0, 34.0, 39.0, -1.0)
    11 = @(0.0, 2.0, 4.0, 0.0, -1.0, -1.0)
    12 = @(0.0, 2.0, 5.0, 244.0, -1.0, 2.0, 38.0, 8.0, 50.0, 16.0, 15.0, 30.0, 20.0, 45.0, 27.0, 14.0, 28.0, 35.0, 36.0, 43.0, -1.0)
}

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y")

foreach ($r in $data.Keys) {
    $values = $data[$r]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $cellRef = "$($cols[$i])$r"
        $ws.Range($cellRef).Value = $values[$i]
    }
}
